$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.93308699131012
$ws.Range("B1").Value = 2.709007740020752
$ws.Range("C1").Value = 2.957734107971191
$ws.Range("D1").Value = 2.732187271118164
$ws.Range("E1").Value = 0.9815295934677124
